$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bots")

# New value in C5, matching the thick-left-border look used by the
# neighboring cells in this column block (e.g. E5 / B5).
$ws.Range("C5").Value = "harry styles"
$ws.Range("C5").Borders.Item(7).Weight = 4

# New row 6 with E6 carrying the same styling.
$ws.Range("E6").Value = "harry styles"
$ws.Range("E6").Borders.Item(7).Weight = 4
